$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 8697
$ws.Range("E2").Value = 530
$ws.Range("F2").Value = 530
$ws.Range("G2").Value = 503
$ws.Range("H2").Value = 368
$ws.Range("I2").Value = 368
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 9004
$ws.Range("L2").Value = 2757
$ws.Range("M2").Value = 6247
$ws.Range("N2").Value = 6222
$ws.Range("O2").Value = 25
$ws.Range("P2").Value = 84
$ws.Range("Q2").Value = 716
$ws.Range("R2").Value = -610
$ws.Range("S2").Value = -239
$ws.Range("T2").Value = 235
$ws.Range("U2").Value = 481
$ws.Range("V2").Value = 1618
$ws.Range("W2").Value = 6.09
$ws.Range("X2").Value = 4.24
$ws.Range("Y2").Value = 6
$ws.Range("Z2").Value = 4.1
$ws.Range("AA2").Value = 44.13
$ws.Range("AB2").Value = 6834.82
$ws.Range("AC2").Value = 21792
$ws.Range("AD2").Value = 7.14
$ws.Range("AE2").Value = 378784
$ws.Range("AF2").Value = 0.41
$ws.Range("AG2").Value = 2500
$ws.Range("AH2").Value = 1.61
$ws.Range("AI2").Value = 11.15
$ws.Range("AJ2").Value = 1690000

# Row 3
$ws.Range("D3").Value = 8258
$ws.Range("E3").Value = 473
$ws.Range("F3").Value = 473
$ws.Range("G3").Value = 467
$ws.Range("H3").Value = 343
$ws.Range("I3").Value = 341
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 8954
$ws.Range("L3").Value = 2385
$ws.Range("M3").Value = 6569
$ws.Range("N3").Value = 6543
$ws.Range("O3").Value = 26
$ws.Range("P3").Value = 84
$ws.Range("Q3").Value = 536
$ws.Range("R3").Value = -467
$ws.Range("S3").Value = -313
$ws.Range("T3").Value = 197
$ws.Range("U3").Value = 339
$ws.Range("V3").Value = 1380
$ws.Range("W3").Value = 5.73
$ws.Range("X3").Value = 4.15
$ws.Range("Y3").Value = 5.35
$ws.Range("Z3").Value = 3.82
$ws.Range("AA3").Value = 36.31
$ws.Range("AB3").Value = 7160.73
$ws.Range("AC3").Value = 20207
$ws.Range("AD3").Value = 9.109999999999999
$ws.Range("AE3").Value = 398323
$ws.Range("AF3").Value = 0.46
$ws.Range("AG3").Value = 2000
$ws.Range("AH3").Value = 1.09
$ws.Range("AI3").Value = 9.619999999999999
$ws.Range("AJ3").Value = 1690000

# Row 4
$ws.Range("D4").Value = 8075
$ws.Range("E4").Value = 442
$ws.Range("F4").Value = 442
$ws.Range("G4").Value = 389
$ws.Range("H4").Value = 274
$ws.Range("I4").Value = 273
$ws.Range("J4").Value = 2
$ws.Range("K4").Value = 8675
$ws.Range("L4").Value = 2056
$ws.Range("M4").Value = 6619
$ws.Range("N4").Value = 6619
$ws.Range("O4").ClearContents()
$ws.Range("P4").Value = 84
$ws.Range("Q4").Value = 862
$ws.Range("R4").Value = -298
$ws.Range("S4").Value = -453
$ws.Range("T4").Value = 169
$ws.Range("U4").Value = 692
$ws.Range("V4").Value = 1030
$ws.Range("W4").Value = 5.47
$ws.Range("X4").Value = 3.4
$ws.Range("Y4").Value = 4.14
$ws.Range("Z4").Value = 3.11
$ws.Range("AA4").Value = 31.06
$ws.Range("AB4").Value = 7441.37
$ws.Range("AC4").Value = 16130
$ws.Range("AD4").Value = 11.16
$ws.Range("AE4").Value = 402922
$ws.Range("AF4").Value = 0.45
$ws.Range("AG4").Value = 2000
$ws.Range("AH4").Value = 1.11
$ws.Range("AI4").Value = 12.05
$ws.Range("AJ4").Value = 1690000

# Row 5
$ws.Range("D5").Value = 8108
$ws.Range("E5").Value = 362
$ws.Range("F5").Value = 362
$ws.Range("G5").Value = 531
$ws.Range("H5").Value = 510
$ws.Range("I5").Value = 510
$ws.Range("J5").ClearContents()
$ws.Range("K5").Value = 8735
$ws.Range("L5").Value = 1631
$ws.Range("M5").Value = 7104
$ws.Range("N5").Value = 7104
$ws.Range("O5").ClearContents()
$ws.Range("P5").Value = 84
$ws.Range("Q5").Value = 434
$ws.Range("R5").Value = -352
$ws.Range("S5").Value = -499
$ws.Range("T5").Value = 158
$ws.Range("U5").Value = 275
$ws.Range("V5").Value = 567
$ws.Range("W5").Value = 4.46
$ws.Range("X5").Value = 6.3
$ws.Range("Y5").Value = 7.44
$ws.Range("Z5").Value = 5.86
$ws.Range("AA5").Value = 22.95
$ws.Range("AB5").Value = 8005.23
$ws.Range("AC5").Value = 30206
$ws.Range("AD5").Value = 5.51
$ws.Range("AE5").Value = 431800
$ws.Range("AF5").Value = 0.39
$ws.Range("AG5").Value = 2000
$ws.Range("AH5").Value = 1.2
$ws.Range("AI5").Value = 6.45
$ws.Range("AJ5").Value = 1690000

# Row 6
$ws.Range("D6").Value = 8646
$ws.Range("E6").Value = 328
$ws.Range("F6").Value = 328
$ws.Range("G6").Value = 749
$ws.Range("H6").Value = 515
$ws.Range("I6").Value = 515
$ws.Range("K6").Value = 9197
$ws.Range("L6").Value = 1707
$ws.Range("M6").Value = 7490
$ws.Range("N6").Value = 7490
$ws.Range("P6").Value = 84
$ws.Range("Q6").Value = 126
$ws.Range("R6").Value = 189
$ws.Range("S6").Value = 77
$ws.Range("T6").Value = 247
$ws.Range("U6").Value = -121
$ws.Range("V6").Value = 567
$ws.Range("W6").Value = 3.79
$ws.Range("X6").Value = 5.95
$ws.Range("Y6").Value = 7.05
$ws.Range("Z6").Value = 5.74
$ws.Range("AA6").Value = 22.79
$ws.Range("AB6").Value = 8578.23
$ws.Range("AC6").Value = 30457
$ws.Range("AD6").Value = 5.78
$ws.Range("AE6").Value = 455260
$ws.Range("AF6").Value = 0.39
$ws.Range("AG6").Value = 2500
$ws.Range("AH6").Value = 1.42
$ws.Range("AI6").Value = 7.99
$ws.Range("AJ6").Value = 1690000

# Rows 7-9: clear financial columns, keep only A/B/C (id, period, company)
$ws.Range("D7:AJ9").ClearContents()

Write-Host "applied edits"